$d = $word.ActiveDocument

# --- Edit 1: "Maybe PMOD joystick." -> "PMOD joystick." -------------------
# The sentence loses its leading "Maybe " word. This paragraph used to hold
# the document's (unique) "_GoBack" bookmark from an earlier edit session;
# since the cursor's last-edit spot moves elsewhere (see Edit 2), drop it
# from here first.
$d.Bookmarks("_GoBack").Delete()

$d.Content.Find.Execute("Maybe PMOD joystick.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "PMOD joystick.", 2) | Out-Null

# --- Edit 2: "Maybe single LEDs." -> "Single LEDs." ------------------------
# Drop "Maybe " and capitalize the leading "s" -> "S".
$d.Content.Find.Execute("Maybe single LEDs.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Single LEDs.", 2) | Out-Null

# Re-seat the "_GoBack" bookmark right after the new "S" (this is now the
# spot of the most recent edit), splitting that run into "S" | "ingle LEDs.".
$hit = $d.Content
$hit.Find.Execute("Single LEDs.", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($hit.Start + 1, $hit.Start + 1)) | Out-Null

# --- Edit 3: footer page field's cached result "1" -> "2" ------------------
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("1", $true, $true, $false, $false, $false, `
                            $true, 1, $false, "2", 2) | Out-Null
